# The commit swaps the OOXML content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml (used by the notes master) had the
# "Office Theme" colors and theme2.xml (used by the slide master, i.e.
# the theme that actually drives the visible deck) had the "Integral"
# colors; after the edit theme2.xml carries the "Office Theme" colors
# (and theme1.xml would carry "Integral", but that part is only wired
# to the notes master and isn't reachable through this COM surface).
#
# Font scheme (Arial/Arial) and format scheme (fills/lines/effects) are
# byte-identical between the two theme parts already, so the only
# observable difference after swapping is the 12-slot color scheme
# that is applied to the slide master's theme. Reproduce that here by
# writing the target ("Office Theme") RGB values into the slide
# master's ThemeColorScheme, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# RGB values packed the usual COLORREF way (R + G*256 + B*65536), since
# the COM-interop runtime here has no RGB() helper function.
$tcs.Item(1).RGB  = 0            # dk1      000000
$tcs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388      # dk2      44546A
$tcs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501      # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407        # accent4  FFC000
$tcs.Item(9).RGB  = 12874308     # accent5  4472C4
$tcs.Item(10).RGB = 4697456      # accent6  70AD47
$tcs.Item(11).RGB = 12673797     # hlink    0563C1
$tcs.Item(12).RGB = 7491477      # folHlink 954F72

# Best-effort: also try to rename the theme / color-scheme back to the
# "Office" naming (a no-op on engines that don't persist this, but
# harmless either way).
try { $master.Theme.Name = "Office Theme" } catch {}
